$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "GW63QV"
$ws.Range("B17").Value = "2024-11-22 21:57:34"
$ws.Range("C17").Value = "GET /"
$ws.Range("D17").Value = 401
$ws.Range("E17").Value = $false
$ws.Range("F17").Value = "Eitss... mau ngapain? Akses terbatas!"

$ws.Range("A18").Value = "V68TJO"
$ws.Range("B18").Value = "2024-11-22 21:57:35"
$ws.Range("C18").Value = "GET /favicon.ico"
$ws.Range("D18").Value = 404
$ws.Range("E18").Value = $false
$ws.Range("F18").Value = "Not Found"
